$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.922.27'
$ws.Range('E2').Value = '  +0.01%  '
$ws.Range('D3').Value = '1.629.85'
$ws.Range('E3').Value = '  -0.55%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '211.89'
$ws.Range('E5').Value = '  -0.77%  '
$ws.Range('E6').Value = '  -0.14%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '23.33'
$ws.Range('E8').Value = '  -1.22%  '
$ws.Range('E9').Value = '  -1.74%  '
$ws.Range('E10').Value = '  -0.27%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0881'
$ws.Range('E11').Value = '  +0.73%  '
$ws.Range('D12').Value = '1.862.05'
$ws.Range('E12').Value = '  -0.49%  '
$ws.Range('D13').Value = '1.627.40'
$ws.Range('E13').Value = '  -0.66%  '
$ws.Range('E14').Value = '  -1.49%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.561'
$ws.Range('E15').Value = '  -2.35%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.61'
$ws.Range('D17').Value = '27.919.04'
$ws.Range('E17').Value = '  +0.04%  '
$ws.Range('E18').Value = '  -0.64%  '
$ws.Range('E19').Value = '  -0.22%  '
$ws.Range('E20').Value = '  +0.73%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.998'
$ws.Range('E21').Value = '  -0.17%  '
$ws.Range('E22').Value = '  -0.46%  '
$ws.Range('E23').Value = '  -5.64%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.04'
$ws.Range('E24').Value = '  -1.13%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '154.70'
$ws.Range('E25').Value = '  +1.95%  '
$ws.Range('E26').Value = '  +0.25%  '
$ws.Range('E27').Value = '  -0.21%  '
$ws.Range('E28').Value = '  -1.23%  '
$ws.Range('E29').Value = '  +0.00%  '
$ws.Range('E30').Value = '  -0.44%  '
$ws.Range('E31').Value = '  -0.37%  '
$ws.Range('E32').Value = '  +2.02%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.07'
$ws.Range('E33').Value = '  -0.82%  '
$ws.Range('D34').Value = '1.398.84'
$ws.Range('E34').Value = '  -0.89%  '
$ws.Range('E35').Value = '  +0.24%  '
$ws.Range('E36').Value = '  +11.26%  '
$ws.Range('E37').Value = '  -0.26%  '
$ws.Range('E38').Value = '  +1.99%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.556'
$ws.Range('E39').Value = '  +0.19%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.863'
$ws.Range('E40').Value = '  -3.17%  '
$ws.Range('E41').Value = '  +0.10%  '
$ws.Range('E42').Value = '  -0.05%  '
$ws.Range('B43').Value = 'RenderToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.84'
$ws.Range('E43').Value = '  +0.31%  '
$ws.Range('B44').Value = 'Aave'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '66.33'
$ws.Range('E44').Value = '  +0.16%  '
$ws.Range('E45').Value = '  +1.11%  '
$ws.Range('E46').Value = '  -0.54%  '
$ws.Range('D47').Value = '1.772.29'
$ws.Range('E47').Value = '  -0.41%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '88.05'
$ws.Range('E49').Value = '  -2.52%  '
$ws.Range('E51').Value = '  -0.09%  '
